$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(45240.451724537037, 0.189500560923064,    "varInd(6)", 5, 10, 0.2, 6, 3, 3),
    @(45240.458680555559, 0.18807162723138701,  "varInd(6)", 5, 10, 0.2, 6, 3, 3),
    @(45240.468680555554, 0.19220388202041699,  "varInd(6)", 5, 10, 0.2, 6, 3, 3)
)

$startRow = 17
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $values = $data[$i]
    $ws.Cells.Item($row, 1).Value = $values[0]
    $ws.Cells.Item($row, 1).NumberFormat = "m/d/yy h:mm"
    $ws.Cells.Item($row, 2).Value = $values[1]
    $ws.Cells.Item($row, 3).Value = $values[2]
    $ws.Cells.Item($row, 4).Value = $values[3]
    $ws.Cells.Item($row, 5).Value = $values[4]
    $ws.Cells.Item($row, 6).Value = $values[5]
    $ws.Cells.Item($row, 7).Value = $values[6]
    $ws.Cells.Item($row, 8).Value = $values[7]
    $ws.Cells.Item($row, 9).Value = $values[8]
}

$ws.Range("J19").Select()
